# Applies the RPA datasets push 2024-05-09 update to the IPO bookbuilding table (Sheet1).
# New entries inserted: 에스오에스랩 (updated dates) and 라메디텍; oldest entry (민테크) rolled off.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (names, date ranges, price ranges, underwriters) ---
$ws.Range("A5").Value = "에스오에스랩"
$ws.Range("B5").Value = "2024.05.28~06.03"
$ws.Range("C5").Value = "7,500~9,000"
$ws.Range("F5").Value = "한국투자증권"
$ws.Range("A6").Value = "디비금융스팩12호"
$ws.Range("B6").Value = "2024.05.28~05.29"
$ws.Range("C6").Value = "2,000~2,000"
$ws.Range("F6").Value = "DB금융투자"
$ws.Range("A7").Value = "씨어스테크놀로지"
$ws.Range("B7").Value = "2024.05.27~05.31"
$ws.Range("C7").Value = "10,500~14,000"
$ws.Range("F7").Value = "한국투자증권"
$ws.Range("A8").Value = "라메디텍"
$ws.Range("B8").Value = "2024.05.27~05.31"
$ws.Range("C8").Value = "10,400~12,700"
$ws.Range("F8").Value = "대신증권"
$ws.Range("A9").Value = "미래에셋비전스팩5호"
$ws.Range("B9").Value = "2024.05.23~05.24"
$ws.Range("C9").Value = "2,000~2,000"
$ws.Range("F9").Value = "미래에셋증권"
$ws.Range("A10").Value = "그리드위즈"
$ws.Range("B10").Value = "2024.05.23~05.29"
$ws.Range("C10").Value = "34,000~40,000"
$ws.Range("F10").Value = "삼성증권"
$ws.Range("A11").Value = "이노스페이스"
$ws.Range("B11").Value = "2024.05.23~05.29"
$ws.Range("C11").Value = "36,400~45,600"
$ws.Range("F11").Value = "미래에셋증권,신한투자증권"
$ws.Range("A12").Value = "에이치브이엠(구.한국진공야금)"
$ws.Range("B12").Value = "2024.05.22~05.28"
$ws.Range("C12").Value = "11,000~14,200"
$ws.Range("F12").Value = "NH투자증권"
$ws.Range("A13").Value = "하스"
$ws.Range("B13").Value = "2024.05.16~05.22"
$ws.Range("C13").Value = "9,000~12,000"
$ws.Range("F13").Value = "삼성증권"
$ws.Range("A14").Value = "미래에셋비전스팩4호"
$ws.Range("B14").Value = "2024.05.13~05.14"
$ws.Range("C14").Value = "2,000~2,000"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "미래에셋증권"
$ws.Range("A15").Value = "노브랜드"
$ws.Range("B15").Value = "2024.04.30~05.08"
$ws.Range("C15").Value = "8,700~11,000"
$ws.Range("D15").Value = "-"
$ws.Range("F15").Value = "삼성증권"
$ws.Range("A16").Value = "KB스팩28호"
$ws.Range("B16").Value = "2024.04.29~04.30"
$ws.Range("C16").Value = "2,000~2,000"
$ws.Range("F16").Value = "KB증권"
$ws.Range("A17").Value = "아이씨티케이"
$ws.Range("B17").Value = "2024.04.24~04.30"
$ws.Range("C17").Value = "13,000~16,000"
$ws.Range("F17").Value = "NH투자증권"
$ws.Range("A18").Value = "SK증권스팩12호"
$ws.Range("B18").Value = "2024.04.17~04.18"
$ws.Range("C18").Value = "2,000~2,000"
$ws.Range("F18").Value = "SK증권"
$ws.Range("A19").Value = "HD현대마린솔루션(구.HD현대글로벌서비스)(유가)"
$ws.Range("B19").Value = "2024.04.16~04.22"
$ws.Range("C19").Value = "73,300~83,400"
$ws.Range("F19").Value = "KB증권,신한투자증권,하나증권,대신증권,삼성증권"
$ws.Range("A20").Value = "코칩"
$ws.Range("B20").Value = "2024.04.15~04.19"
$ws.Range("C20").Value = "11,000~14,000"
$ws.Range("F20").Value = "한국투자증권"
$ws.Range("A21").Value = "유안타스팩16호"
$ws.Range("B21").Value = "2024.04.15~04.16"
$ws.Range("C21").Value = "2,000~2,000"
$ws.Range("F21").Value = "유안타증권"

# --- Numeric-looking cells that must remain stored as text (matches source formatting) ---
$numericTextCells = @("E5","E6","E7","E8","E9","E10","E11","E12","E13","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21")
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}
$ws.Range("E5").Value = "15000"
$ws.Range("E6").Value = "10000"
$ws.Range("E7").Value = "13650"
$ws.Range("E8").Value = "13499"
$ws.Range("E9").Value = "9500"
$ws.Range("E10").Value = "47600"
$ws.Range("E11").Value = "48412"
$ws.Range("E12").Value = "26400"
$ws.Range("E13").Value = "16290"
$ws.Range("E15").Value = "10440"
$ws.Range("D16").Value = "2000"
$ws.Range("E16").Value = "10000"
$ws.Range("D17").Value = "20000"
$ws.Range("E17").Value = "25610"
$ws.Range("D18").Value = "2000"
$ws.Range("E18").Value = "6000"
$ws.Range("D19").Value = "83400"
$ws.Range("E19").Value = "652370"
$ws.Range("D20").Value = "18000"
$ws.Range("E20").Value = "16500"
$ws.Range("D21").Value = "2000"
$ws.Range("E21").Value = "10300"
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).Style = "Normal"
}
